$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-31 Friday" "2024-06-01 Saturday"

Replace-Text "25×34=850" "83×12=996"
Replace-Text "64×30=1920" "97×92=8924"
Replace-Text "97×11=1067" "54×71=3834"
Replace-Text "65×13=845" "18×36=648"
Replace-Text "51×37=1887" "34×43=1462"

Replace-Text "56×24=1344" "93×18=1674"
Replace-Text "19×87=1653" "60×27=1620"
Replace-Text "82×16=1312" "29×49=1421"
Replace-Text "61×61=3721" "22×35=770"
Replace-Text "98×58=5684" "25×77=1925"

Replace-Text "90×32=2880" "70×51=3570"
Replace-Text "79×49=3871" "83×70=5810"
Replace-Text "88×92=8096" "56×52=2912"
Replace-Text "54×95=5130" "72×84=6048"
Replace-Text "15×17=255" "99×92=9108"

Replace-Text "27×74=1998" "63×17=1071"
Replace-Text "34×74=2516" "19×72=1368"
Replace-Text "74×60=4440" "86×50=4300"
Replace-Text "72×27=1944" "73×39=2847"
Replace-Text "51×35=1785" "28×36=1008"

Replace-Text "74×61=4514" "67×28=1876"
Replace-Text "83×95=7885" "58×82=4756"
Replace-Text "76×25=1900" "92×33=3036"
Replace-Text "47×23=1081" "78×38=2964"
Replace-Text "41×46=1886" "97×68=6596"
